$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44571
$ws.Cells.Item(2, 11).Value = 'Brooks'
$ws.Cells.Item(2, 12).Value = 'Segunda'
$ws.Cells.Item(2, 13).Value = 400
$ws.Cells.Item(2, 14).Value = 8500
$ws.Cells.Item(2, 15).Value = 9000
$ws.Cells.Item(2, 16).Value = 8750
$ws.Cells.Item(2, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 875
$ws.Cells.Item(2, 20).Value = 10

# Row 3
$ws.Cells.Item(3, 4).Value = 44208
$ws.Cells.Item(3, 11).Value = 'Lapins'
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 10500
$ws.Cells.Item(3, 15).Value = 11000
$ws.Cells.Item(3, 16).Value = 10750
$ws.Cells.Item(3, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value = 896
$ws.Cells.Item(3, 20).Value = 12

# Row 5
$ws.Cells.Item(5, 4).Value = 44568
$ws.Cells.Item(5, 11).Value = 'Santina'
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 15000
$ws.Cells.Item(5, 15).Value = 16000
$ws.Cells.Item(5, 16).Value = 15500
$ws.Cells.Item(5, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(5, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 19).Value = 1292
$ws.Cells.Item(5, 20).Value = 12

# Row 6
$ws.Cells.Item(6, 4).Value = 44580
$ws.Cells.Item(6, 11).Value = 'Sweet Heart'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 7000
$ws.Cells.Item(6, 15).Value = 8000
$ws.Cells.Item(6, 16).Value = 7500
$ws.Cells.Item(6, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(6, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 19).Value = 750
$ws.Cells.Item(6, 20).Value = 10

# Row 8
$ws.Cells.Item(8, 4).Value = 44537
$ws.Cells.Item(8, 11).Value = 'Brooks'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 29000
$ws.Cells.Item(8, 15).Value = 30000
$ws.Cells.Item(8, 16).Value = 29500
$ws.Cells.Item(8, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1475
$ws.Cells.Item(8, 20).Value = 20

# Row 10
$ws.Cells.Item(10, 4).Value = 44532
$ws.Cells.Item(10, 11).Value = 'Brooks'
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 27000
$ws.Cells.Item(10, 15).Value = 28000
$ws.Cells.Item(10, 16).Value = 27500
$ws.Cells.Item(10, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 2292
$ws.Cells.Item(10, 20).Value = 12

# Row 11
$ws.Cells.Item(11, 4).Value = 44594
$ws.Cells.Item(11, 11).Value = 'Santina'
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 160
$ws.Cells.Item(11, 14).Value = 5000
$ws.Cells.Item(11, 15).Value = 6000
$ws.Cells.Item(11, 16).Value = 5500
$ws.Cells.Item(11, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(11, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 19).Value = 1100
$ws.Cells.Item(11, 20).Value = 5

# Row 12
$ws.Cells.Item(12, 4).Value = 44210
$ws.Cells.Item(12, 11).Value = 'Rainier'
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 13).Value = 250
$ws.Cells.Item(12, 14).Value = 21000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 21500
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 1194
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = 44557
$ws.Cells.Item(13, 11).Value = 'Lapins'
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 250
$ws.Cells.Item(13, 14).Value = 9000
$ws.Cells.Item(13, 15).Value = 10000
$ws.Cells.Item(13, 16).Value = 9500
$ws.Cells.Item(13, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(13, 19).Value = 950
$ws.Cells.Item(13, 20).Value = 10
